$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-09-16 to 2023-10-05
# (Excel serial date 45204 == 2023-10-05)
$newDateSerial = 45204

$ws.Range("C2").Value = $newDateSerial
$ws.Range("C3").Value = $newDateSerial
$ws.Range("C4").Value = $newDateSerial
$ws.Range("C5").Value = $newDateSerial
